# B6-PowerPoint.pptx edit
#
# 1) The three tables (on slides 14, 15 and 16) switch from the bespoke
#    "Table_0" style ({F521FE54-7733-4165-BDC3-3226D920392D}) to the
#    built-in "No Style, Table Grid" style
#    ({1011830E-FE34-4C89-A9B9-B9B218F9172A}).
#
# 2) The deck's theme (ppt/theme/theme1.xml, used by the slide master /
#    all slides) is switched from the "Integral" / "Red Violet" palette
#    over to the stock "Office" palette (the palette that used to live,
#    unused, in ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$tableSlides = 14, 15, 16
$newTableStyle = "{1011830E-FE34-4C89-A9B9-B9B218F9172A}"

foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the theme color palette --------------------------------------
# New (target) "Office" scheme values, in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order, expressed as the
# decimal BGR-packed long that the PowerPoint object model's
# ColorFormat/ThemeColor.RGB property expects.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
